$wb = $excel.ActiveWorkbook

$wsInfo  = $wb.Worksheets.Item("info")
$wsGWP   = $wb.Worksheets.Item("GWP")
$wsEutro = $wb.Worksheets.Item("Eutrophication")

# Rename "Bricks" -> "Brick" on the info sheet (row 3).
$wsInfo.Range("A3").Value = "Brick"

# GWP and Eutrophication: point column A (rows 2-12) at the info sheet via
# formula instead of a static label, so the material names stay linked.
for ($r = 2; $r -le 12; $r++) {
    $wsGWP.Range("A$r").Formula   = "=info!A$r"
    $wsEutro.Range("A$r").Formula = "=info!A$r"
}

# Restore/update each sheet's selection.
$wsInfo.Range("A4").Select()
$wsGWP.Range("A2:A12").Select()

# Eutrophication becomes the active (last-viewed) sheet/tab.
$wsEutro.Activate()
$wsEutro.Range("B16").Select()
